# Sheet1 is a flat "label / metric / unit / year / value" extraction table.
# The commit relabels several metric blocks (energy/GHG/water/waste got
# shuffled one slot down within their groups), fixes a unit typo
# (MWh -> MWhs), and appends 5 rows that were missing from the original
# extraction (water-consumption-intensity 2021 + waste-generated 2023/2022),
# growing the used range from A1:F33 to A1:F38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (metric) relabels: rows 5-13 all become "Total energy consumption" ---
$ws.Cells.Item(5, 3).Value  = "Total energy consumption"
$ws.Cells.Item(6, 3).Value  = "Total energy consumption"
$ws.Cells.Item(7, 3).Value  = "Total energy consumption"
$ws.Cells.Item(8, 3).Value  = "Total energy consumption"
$ws.Cells.Item(9, 3).Value  = "Total energy consumption"
$ws.Cells.Item(10, 3).Value = "Total energy consumption"
$ws.Cells.Item(11, 3).Value = "Total energy consumption"
$ws.Cells.Item(11, 4).Value = "MWhs"
$ws.Cells.Item(12, 3).Value = "Total energy consumption"
$ws.Cells.Item(12, 4).Value = "MWhs"
$ws.Cells.Item(13, 3).Value = "Total energy consumption"
$ws.Cells.Item(13, 4).Value = "MWhs"

# --- Rows 14-16: was "Greenhouse Gas Emissions / Total / tCO2e", becomes the
#     new "Energy Consumption / Energy consumption intensity / MWhs" block ---
$ws.Cells.Item(14, 2).Value = "Energy Consumption"
$ws.Cells.Item(14, 3).Value = "Energy consumption intensity"
$ws.Cells.Item(14, 4).Value = "MWhs"
$ws.Cells.Item(14, 6).Value = 0.021
$ws.Cells.Item(15, 2).Value = "Energy Consumption"
$ws.Cells.Item(15, 3).Value = "Energy consumption intensity"
$ws.Cells.Item(15, 4).Value = "MWhs"
$ws.Cells.Item(15, 6).Value = 0.02
$ws.Cells.Item(16, 2).Value = "Energy Consumption"
$ws.Cells.Item(16, 3).Value = "Energy consumption intensity"
$ws.Cells.Item(16, 4).Value = "MWhs"
$ws.Cells.Item(16, 6).Value = 0.022

# --- Rows 17-25: GHG metric labels + values all shift down one slot
#     (Scope1->Total, Scope2->Scope1, Scope3->Scope2), each keeping the
#     value previously held by the row above it ---
$ws.Cells.Item(17, 3).Value = "Total"
$ws.Cells.Item(17, 6).Value = 27497
$ws.Cells.Item(18, 3).Value = "Total"
$ws.Cells.Item(18, 6).Value = 65488
$ws.Cells.Item(19, 3).Value = "Total"
$ws.Cells.Item(19, 6).Value = 68151
$ws.Cells.Item(20, 3).Value = "Scope 1"
$ws.Cells.Item(20, 6).Value = 147
$ws.Cells.Item(21, 3).Value = "Scope 1"
$ws.Cells.Item(21, 6).Value = 160
$ws.Cells.Item(22, 3).Value = "Scope 1"
$ws.Cells.Item(22, 6).Value = 237
$ws.Cells.Item(23, 3).Value = "Scope 2"
$ws.Cells.Item(23, 6).Value = 68334
$ws.Cells.Item(24, 3).Value = "Scope 2"
$ws.Cells.Item(24, 6).Value = 63811
$ws.Cells.Item(25, 3).Value = "Scope 2"
$ws.Cells.Item(25, 6).Value = 67636

# --- Rows 26-28: was "Water Consumption / Total water consumption / m³",
#     becomes "Greenhouse Gas Emissions / Scope 3 / tCO2e" (continuing the
#     GHG block's downward shift), keeping the GHG values ---
$ws.Cells.Item(26, 2).Value = "Greenhouse Gas Emissions"
$ws.Cells.Item(26, 3).Value = "Scope 3"
$ws.Cells.Item(26, 4).Value = "tCO2e"
$ws.Cells.Item(26, 6).Value = 3849
$ws.Cells.Item(27, 2).Value = "Greenhouse Gas Emissions"
$ws.Cells.Item(27, 3).Value = "Scope 3"
$ws.Cells.Item(27, 4).Value = "tCO2e"
$ws.Cells.Item(27, 6).Value = 1517
$ws.Cells.Item(28, 2).Value = "Greenhouse Gas Emissions"
$ws.Cells.Item(28, 3).Value = "Scope 3"
$ws.Cells.Item(28, 4).Value = "tCO2e"
$ws.Cells.Item(28, 6).Value = 278

# --- Rows 29-30: was "Water Consumption / Water consumption intensity /
#     m³/ft2", becomes "Greenhouse Gas Emissions / Emission intensities of
#     Scope 2 / tCO2e" with new intensity values ---
$ws.Cells.Item(29, 2).Value = "Greenhouse Gas Emissions"
$ws.Cells.Item(29, 3).Value = "Emission intensities of Scope 2"
$ws.Cells.Item(29, 4).Value = "tCO2e"
$ws.Cells.Item(29, 6).Value = 0.0123
$ws.Cells.Item(30, 2).Value = "Greenhouse Gas Emissions"
$ws.Cells.Item(30, 3).Value = "Emission intensities of Scope 2"
$ws.Cells.Item(30, 4).Value = "tCO2e"
$ws.Cells.Item(30, 6).Value = 0.0115

# --- Rows 31-33: was "Waste Generation / Total waste generated / t",
#     becomes "Water Consumption / Total water consumption / m³", taking on
#     the values that used to live in rows 26-28 (row 33 is a brand-new
#     water-consumption data point, value 407051, replacing the old
#     "waste generated 2021 = 379" row) ---
$ws.Cells.Item(31, 2).Value = "Water Consumption"
$ws.Cells.Item(31, 3).Value = "Total water consumption"
$ws.Cells.Item(31, 4).Value = "m³"
$ws.Cells.Item(31, 6).Value = 433969
$ws.Cells.Item(32, 2).Value = "Water Consumption"
$ws.Cells.Item(32, 3).Value = "Total water consumption"
$ws.Cells.Item(32, 4).Value = "m³"
$ws.Cells.Item(32, 6).Value = 400322
$ws.Cells.Item(33, 2).Value = "Water Consumption"
$ws.Cells.Item(33, 3).Value = "Total water consumption"
$ws.Cells.Item(33, 4).Value = "m³"
$ws.Cells.Item(33, 6).Value = 407051

# --- Brand-new rows 34-38 (the sheet grows from A1:F33 to A1:F38) ---

# Column A on the data rows carries a bordered/centered style (s="1"); copy
# that formatting from the last existing data row before filling A34:A38 so
# the new index cells match the rest of the column.
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A34:A38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The "year" column holds text like "2023", not numbers - format the new
# cells as Text first so they don't get auto-coerced to numeric.
$ws.Range("E34:E38").NumberFormat = "@"

$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = "Water Consumption"
$ws.Cells.Item(34, 3).Value = "Water consumption intensity"
$ws.Cells.Item(34, 4).Value = "m³"
$ws.Cells.Item(34, 5).Value = "2023"
$ws.Cells.Item(34, 6).Value = 0.101

$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = "Water Consumption"
$ws.Cells.Item(35, 3).Value = "Water consumption intensity"
$ws.Cells.Item(35, 4).Value = "m³"
$ws.Cells.Item(35, 5).Value = "2022"
$ws.Cells.Item(35, 6).Value = 0.097

$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = "Water Consumption"
$ws.Cells.Item(36, 3).Value = "Water consumption intensity"
$ws.Cells.Item(36, 4).Value = "m³"
$ws.Cells.Item(36, 5).Value = "2021"
$ws.Cells.Item(36, 6).Value = 0.077

$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = "Waste Generation"
$ws.Cells.Item(37, 3).Value = "Total waste generated"
$ws.Cells.Item(37, 4).Value = "t"
$ws.Cells.Item(37, 5).Value = "2023"
$ws.Cells.Item(37, 6).Value = 1605

$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = "Waste Generation"
$ws.Cells.Item(38, 3).Value = "Total waste generated"
$ws.Cells.Item(38, 4).Value = "t"
$ws.Cells.Item(38, 5).Value = "2022"
$ws.Cells.Item(38, 6).Value = 1229
